$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 937.5476
$ws.Range("J17").Value = 937.5476
$ws.Range("L17").Value = 2812.6428
$ws.Range("N17").Value = -3148.6428
$ws.Range("H28").Value = 187.54546
$ws.Range("I28").Value = 181.75
$ws.Range("J28").Value = 203
$ws.Range("K28").Value = 181.75
$ws.Range("L28").Value = 203
$ws.Range("M28").Value = 303.25
$ws.Range("N28").Value = -1173
$ws.Range("H40").Value = 2299.1428
$ws.Range("I40").Value = 3824.5
$ws.Range("J40").Value = 1689
$ws.Range("K40").Value = 3824.5
$ws.Range("L40").Value = 1689
$ws.Range("M40").Value = -3649.5
$ws.Range("N40").Value = -2039
$ws.Range("H76").Value = 6081.091
$ws.Range("I76").Value = 5421.2
$ws.Range("J76").Value = 6631
$ws.Range("K76").Value = 5421.2
$ws.Range("L76").Value = 6631
$ws.Range("M76").Value = -5106.2
$ws.Range("N76").Value = -7261
$ws.Range("H79").Value = 6081.091
$ws.Range("I79").Value = 5421.2
$ws.Range("J79").Value = 6631
$ws.Range("K79").Value = 5421.2
$ws.Range("L79").Value = 6631
$ws.Range("M79").Value = -4329.2
$ws.Range("N79").Value = -8815
$ws.Range("H92").Value = 816.82355
$ws.Range("I92").Value = 695.73334
$ws.Range("K92").Value = 695.73334
$ws.Range("M92").Value = 552.26666
$ws.Range("H132").Value = 6948750
$ws.Range("I132").Value = 7411740
$ws.Range("K132").Value = 22235220
$ws.Range("M132").Value = -22232690
$ws.Range("H141").Value = 3000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 9000
$ws.Range("N141").Value = -19360
$ws.Range("M141").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1112.6957
$ws.Range("I61").Value = 849.5
$ws.Range("J61").Value = 1714.2858
$ws.Range("K61").Value = 849.5
$ws.Range("L61").Value = 1714.2858
$ws.Range("M61").Value = -637.5
$ws.Range("N61").Value = -2138.2858
$ws.Range("H74").Value = 1230.6471
$ws.Range("I74").Value = 552.3043
$ws.Range("K74").Value = 552.3043
$ws.Range("M74").Value = 321.6957
$ws.Range("H77").Value = 1230.6471
$ws.Range("I77").Value = 552.3043
$ws.Range("K77").Value = 2761.5215
$ws.Range("M77").Value = 1606.4785
$ws.Range("H112").Value = 29593.5
$ws.Range("J112").Value = 29593.5
$ws.Range("L112").Value = 29593.5
$ws.Range("N112").Value = -32547.5
$ws.Range("H132").Value = 1943.65
$ws.Range("I132").Value = 1675.1666
$ws.Range("J132").Value = 2749.1
$ws.Range("K132").Value = 5025.4998
$ws.Range("L132").Value = 8247.299999999999
$ws.Range("M132").Value = -2495.4998
$ws.Range("N132").Value = -13307.3
$ws.Range("H136").Value = 1112.6957
$ws.Range("I136").Value = 849.5
$ws.Range("J136").Value = 1714.2858
$ws.Range("K136").Value = 2548.5
$ws.Range("L136").Value = 5142.857400000001
$ws.Range("M136").Value = 1.5
$ws.Range("N136").Value = -10242.8574
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3182.5454
$ws.Range("I86").Value = 3305.6667
$ws.Range("J86").Value = 2967.0833
$ws.Range("K86").Value = 3305.6667
$ws.Range("L86").Value = 2967.0833
$ws.Range("M86").Value = -2182.6667
$ws.Range("N86").Value = -5213.0833
$ws.Range("H89").Value = 3182.5454
$ws.Range("I89").Value = 3305.6667
$ws.Range("J89").Value = 2967.0833
$ws.Range("K89").Value = 16528.3335
$ws.Range("L89").Value = 14835.4165
$ws.Range("M89").Value = -10912.3335
$ws.Range("N89").Value = -26067.4165
$ws.Range("H110").Value = 40427.832
$ws.Range("J110").Value = 40427.832
$ws.Range("L110").Value = 40427.832
$ws.Range("N110").Value = -48607.832
$ws.Range("H138").Value = 67186.664
$ws.Range("J138").Value = 67186.664
$ws.Range("L138").Value = 67186.664
$ws.Range("N138").Value = -77466.664
$ws.Range("H140").Value = 18000
$ws.Range("J140").Value = 18000
$ws.Range("L140").Value = 18000
$ws.Range("N140").Value = -28360
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 221040.2
$ws.Range("I3").Value = 49600.5
$ws.Range("J3").Value = 335333.34
$ws.Range("K3").Value = 49600.5
$ws.Range("L3").Value = 335333.34
$ws.Range("M3").Value = -49487.5
$ws.Range("N3").Value = -335559.34
$ws.Range("H58").Value = 1362.8823
$ws.Range("I58").Value = 1150.909
$ws.Range("J58").Value = 1751.5
$ws.Range("K58").Value = 1150.909
$ws.Range("L58").Value = 1751.5
$ws.Range("M58").Value = -947.9090000000001
$ws.Range("N58").Value = -2157.5
$ws.Range("H136").Value = 1362.8823
$ws.Range("I136").Value = 1150.909
$ws.Range("J136").Value = 1751.5
$ws.Range("K136").Value = 3452.727
$ws.Range("L136").Value = 5254.5
$ws.Range("M136").Value = -902.7270000000003
$ws.Range("N136").Value = -10354.5
$ws.Range("H141").Value = 29014.285
$ws.Range("J141").Value = 29014.285
$ws.Range("L141").Value = 29014.285
$ws.Range("N141").Value = -39374.285
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1405.375
$ws.Range("J68").Value = 1544.8182
$ws.Range("L68").Value = 4634.4546
$ws.Range("N68").Value = -6256.4546
$ws.Range("H71").Value = 1405.375
$ws.Range("J71").Value = 1544.8182
$ws.Range("L71").Value = 13903.3638
$ws.Range("N71").Value = -22015.3638
$ws.Range("H121").Value = 860.7273
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 896.8
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 2690.4
$ws.Range("N121").Value = -5310.4
$ws.Range("M121").Value = -190
$ws.Range("H134").Value = 3459
$ws.Range("I134").Value = 1649.9375
$ws.Range("J134").Value = 8283.166999999999
$ws.Range("K134").Value = 4949.8125
$ws.Range("L134").Value = 24849.501
$ws.Range("M134").Value = 120.1875
$ws.Range("N134").Value = -34989.501
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1241.8125
$ws.Range("I122").Value = 1309.2916
$ws.Range("K122").Value = 3927.8748
$ws.Range("M122").Value = -1477.8748
$ws.Range("H132").Value = 1902.9412
$ws.Range("I132").Value = 1700.2
$ws.Range("K132").Value = 5100.6
$ws.Range("M132").Value = -2570.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2602.9644
$ws.Range("I40").Value = 2317.8823
$ws.Range("J40").Value = 3043.5454
$ws.Range("K40").Value = 2317.8823
$ws.Range("L40").Value = 3043.5454
$ws.Range("M40").Value = -2181.8823
$ws.Range("N40").Value = -3315.5454
$ws.Range("H55").Value = 268.65216
$ws.Range("I55").Value = 177.9
$ws.Range("J55").Value = 338.46155
$ws.Range("K55").Value = 177.9
$ws.Range("L55").Value = 338.46155
$ws.Range("M55").Value = -4.900000000000006
$ws.Range("N55").Value = -684.46155
$ws.Range("H110").Value = 30429
$ws.Range("J110").Value = 30643.5
$ws.Range("L110").Value = 30643.5
$ws.Range("N110").Value = -38823.5
$ws.Range("H132").Value = 21297.295
$ws.Range("I132").Value = 1258.1936
$ws.Range("K132").Value = 3774.5808
$ws.Range("M132").Value = -1244.5808
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2102.7942
$ws.Range("I132").Value = 1710.1428
$ws.Range("J132").Value = 2737.077
$ws.Range("K132").Value = 5130.428400000001
$ws.Range("L132").Value = 8211.231
$ws.Range("M132").Value = -2600.428400000001
$ws.Range("N132").Value = -13271.231
